$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Insert 24 new rows before row 5 (shifting the existing rows 5-28 down to 29-52),
# matching the net change introduced by the new debug screens / custom containers.
$ws.Range("A5:A28").EntireRow.Insert()

$ws.Range("B5").Value = 'SingleUseId25'
$ws.Range("C5").Value = 'Default'
$ws.Range("D5").Value = 'Center'
$ws.Range("E5").Value = 'LTR'
$ws.Range("F5").Value = 'Debug'
$ws.Range("B6").Value = 'SingleUseId26'
$ws.Range("C6").Value = 'Default'
$ws.Range("D6").Value = 'Center'
$ws.Range("E6").Value = 'LTR'
$ws.Range("F6").Value = 'Application'
$ws.Range("B7").Value = 'SingleUseId28'
$ws.Range("C7").Value = 'Default'
$ws.Range("D7").Value = 'Center'
$ws.Range("E7").Value = 'LTR'
$ws.Range("F7").Value = 'Application'
$ws.Range("B8").Value = 'SingleUseId29'
$ws.Range("C8").Value = 'Default'
$ws.Range("D8").Value = 'Center'
$ws.Range("E8").Value = 'LTR'
$ws.Range("F8").Value = 'New Text'
$ws.Range("B9").Value = 'SingleUseId31'
$ws.Range("C9").Value = 'Default'
$ws.Range("D9").Value = 'Left'
$ws.Range("E9").Value = 'LTR'
$ws.Range("F9").Value = 'New Text'
$ws.Range("B10").Value = 'SingleUseId32'
$ws.Range("C10").Value = 'Default'
$ws.Range("D10").Value = 'Left'
$ws.Range("E10").Value = 'LTR'
$ws.Range("F10").Value = 'New Text'
$ws.Range("B11").Value = 'SingleUseId33'
$ws.Range("C11").Value = 'Default'
$ws.Range("D11").Value = 'Left'
$ws.Range("E11").Value = 'LTR'
$ws.Range("F11").Value = 'New Text'
$ws.Range("B12").Value = 'SingleUseId34'
$ws.Range("C12").Value = 'Default'
$ws.Range("D12").Value = 'Left'
$ws.Range("E12").Value = 'LTR'
$ws.Range("F12").Value = 'New Text'
$ws.Range("B13").Value = 'SingleUseId35'
$ws.Range("C13").Value = 'Default'
$ws.Range("D13").Value = 'Left'
$ws.Range("E13").Value = 'LTR'
$ws.Range("F13").Value = 'New Text'
$ws.Range("B14").Value = 'SingleUseId36'
$ws.Range("C14").Value = 'Default'
$ws.Range("D14").Value = 'Left'
$ws.Range("E14").Value = 'LTR'
$ws.Range("F14").Value = 'New Text'
$ws.Range("B15").Value = 'SingleUseId37'
$ws.Range("C15").Value = 'Default'
$ws.Range("D15").Value = 'Left'
$ws.Range("E15").Value = 'LTR'
$ws.Range("F15").Value = 'New Text'
$ws.Range("B16").Value = 'SingleUseId38'
$ws.Range("C16").Value = 'Default'
$ws.Range("D16").Value = 'Left'
$ws.Range("E16").Value = 'LTR'
$ws.Range("F16").Value = 'New Text'
$ws.Range("B17").Value = 'SingleUseId39'
$ws.Range("C17").Value = 'Default'
$ws.Range("D17").Value = 'Left'
$ws.Range("E17").Value = 'LTR'
$ws.Range("F17").Value = 'New Text'
$ws.Range("B18").Value = 'SingleUseId40'
$ws.Range("C18").Value = 'Small'
$ws.Range("D18").Value = 'Left'
$ws.Range("E18").Value = 'LTR'
$ws.Range("F18").Value = 'Time: <value>'
$ws.Range("B19").Value = 'SingleUseId41'
$ws.Range("C19").Value = 'Small'
$ws.Range("D19").Value = 'Left'
$ws.Range("E19").Value = 'LTR'
$ws.Range("F19").Value = '00:00:00'
$ws.Range("B20").Value = 'SingleUseId42'
$ws.Range("C20").Value = 'Small'
$ws.Range("D20").Value = 'Left'
$ws.Range("E20").Value = 'LTR'
$ws.Range("F20").Value = 'Lat: <value> <value>'
$ws.Range("B21").Value = 'SingleUseId43'
$ws.Range("C21").Value = 'Small'
$ws.Range("D21").Value = 'Left'
$ws.Range("E21").Value = 'LTR'
$ws.Range("F21").Value = '---'
$ws.Range("B22").Value = 'SingleUseId44'
$ws.Range("C22").Value = 'Small'
$ws.Range("D22").Value = 'Left'
$ws.Range("E22").Value = 'LTR'
$ws.Range("F22").Value = 'Lon: <value> <value>'
$ws.Range("B23").Value = 'SingleUseId45'
$ws.Range("C23").Value = 'Small'
$ws.Range("D23").Value = 'Left'
$ws.Range("E23").Value = 'LTR'
$ws.Range("F23").Value = '---'
$ws.Range("B24").Value = 'SingleUseId46'
$ws.Range("C24").Value = 'Small'
$ws.Range("D24").Value = 'Left'
$ws.Range("E24").Value = 'LTR'
$ws.Range("F24").Value = 'Alti: <value> m'
$ws.Range("B25").Value = 'SingleUseId47'
$ws.Range("C25").Value = 'Small'
$ws.Range("D25").Value = 'Left'
$ws.Range("E25").Value = 'LTR'
$ws.Range("F25").Value = '---'
$ws.Range("B26").Value = 'SingleUseId48'
$ws.Range("C26").Value = 'Small'
$ws.Range("D26").Value = 'Left'
$ws.Range("E26").Value = 'LTR'
$ws.Range("F26").Value = 'Fix: <value>'
$ws.Range("B27").Value = 'SingleUseId49'
$ws.Range("C27").Value = 'Small'
$ws.Range("D27").Value = 'Left'
$ws.Range("E27").Value = 'LTR'
$ws.Range("F27").Value = '---'
$ws.Range("B28").Value = 'SingleUseId50'
$ws.Range("C28").Value = 'Small'
$ws.Range("D28").Value = 'Left'
$ws.Range("E28").Value = 'LTR'
$ws.Range("F28").Value = 'Sat: <value>'
$ws.Range("B29").Value = 'SingleUseId51'
$ws.Range("C29").Value = 'Small'
$ws.Range("D29").Value = 'Left'
$ws.Range("E29").Value = 'LTR'
$ws.Range("F29").Value = '---'
$ws.Range("B30").Value = 'SingleUseId52'
$ws.Range("C30").Value = 'Small'
$ws.Range("D30").Value = 'Left'
$ws.Range("E30").Value = 'LTR'
$ws.Range("F30").Value = 'Date: <value>'
$ws.Range("B31").Value = 'SingleUseId53'
$ws.Range("C31").Value = 'Small'
$ws.Range("D31").Value = 'Left'
$ws.Range("E31").Value = 'LTR'
$ws.Range("F31").Value = '01.01.2020'
$ws.Range("B32").Value = 'SingleUseId54'
$ws.Range("C32").Value = 'Small'
$ws.Range("D32").Value = 'Left'
$ws.Range("E32").Value = 'LTR'
$ws.Range("F32").Value = 'Read: <value>'
$ws.Range("B33").Value = 'SingleUseId55'
$ws.Range("C33").Value = 'Small'
$ws.Range("D33").Value = 'Left'
$ws.Range("E33").Value = 'LTR'
$ws.Range("F33").Value = '---'
$ws.Range("B34").Value = 'SingleUseId56'
$ws.Range("C34").Value = 'Small'
$ws.Range("D34").Value = 'Left'
$ws.Range("E34").Value = 'LTR'
$ws.Range("F34").Value = 'Write: <value>'
$ws.Range("B35").Value = 'SingleUseId57'
$ws.Range("C35").Value = 'Small'
$ws.Range("D35").Value = 'Left'
$ws.Range("E35").Value = 'LTR'
$ws.Range("F35").Value = '---'
$ws.Range("B36").Value = 'SingleUseId58'
$ws.Range("C36").Value = 'Default'
$ws.Range("D36").Value = 'Center'
$ws.Range("E36").Value = 'LTR'
$ws.Range("F36").Value = 'Debug Gps Data'
$ws.Range("B37").Value = 'SingleUseId60'
$ws.Range("C37").Value = 'Small'
$ws.Range("D37").Value = 'Center'
$ws.Range("E37").Value = 'LTR'
$ws.Range("F37").Value = 'Next'
$ws.Range("B38").Value = 'SingleUseId63'
$ws.Range("C38").Value = 'Default'
$ws.Range("D38").Value = 'Center'
$ws.Range("E38").Value = 'LTR'
$ws.Range("F38").Value = 'Debug Draw'
$ws.Range("B39").Value = 'SingleUseId67'
$ws.Range("C39").Value = 'Small'
$ws.Range("D39").Value = 'Center'
$ws.Range("E39").Value = 'LTR'
$ws.Range("F39").Value = 'Exit'
$ws.Range("B40").Value = 'SingleUseId68'
$ws.Range("C40").Value = 'Small'
$ws.Range("D40").Value = 'Center'
$ws.Range("E40").Value = 'LTR'
$ws.Range("F40").Value = 'Next'
$ws.Range("B41").Value = 'SingleUseId69'
$ws.Range("C41").Value = 'Default'
$ws.Range("D41").Value = 'Center'
$ws.Range("E41").Value = 'LTR'
$ws.Range("F41").Value = 'Debug SD card'
$ws.Range("B42").Value = 'SingleUseId59'
$ws.Range("C42").Value = 'Small'
$ws.Range("D42").Value = 'Center'
$ws.Range("E42").Value = 'LTR'
$ws.Range("F42").Value = 'Next'
$ws.Range("B43").Value = 'SingleUseId70'
$ws.Range("C43").Value = 'Small'
$ws.Range("D43").Value = 'Left'
$ws.Range("E43").Value = 'LTR'
$ws.Range("F43").Value = 'State: <value>'
$ws.Range("B44").Value = 'SingleUseId71'
$ws.Range("C44").Value = 'Small'
$ws.Range("D44").Value = 'Left'
$ws.Range("E44").Value = 'LTR'
$ws.Range("F44").Value = 'UNINITIALIZED'
$ws.Range("B45").Value = 'SingleUseId72'
$ws.Range("C45").Value = 'Small'
$ws.Range("D45").Value = 'Left'
$ws.Range("E45").Value = 'LTR'
$ws.Range("F45").Value = 'Total space: <value>'
$ws.Range("B46").Value = 'SingleUseId73'
$ws.Range("C46").Value = 'Small'
$ws.Range("D46").Value = 'Left'
$ws.Range("E46").Value = 'LTR'
$ws.Range("F46").Value = '0'
$ws.Range("B47").Value = 'SingleUseId74'
$ws.Range("C47").Value = 'Small'
$ws.Range("D47").Value = 'Left'
$ws.Range("E47").Value = 'LTR'
$ws.Range("F47").Value = 'Free space: <value>'
$ws.Range("B48").Value = 'SingleUseId75'
$ws.Range("C48").Value = 'Small'
$ws.Range("D48").Value = 'Left'
$ws.Range("E48").Value = 'LTR'
$ws.Range("F48").Value = '0'
$ws.Range("B49").Value = 'SingleUseId76'
$ws.Range("C49").Value = 'Small'
$ws.Range("D49").Value = 'Left'
$ws.Range("E49").Value = 'LTR'
$ws.Range("F49").Value = 'Input files: <value>'
$ws.Range("B50").Value = 'SingleUseId77'
$ws.Range("C50").Value = 'Small'
$ws.Range("D50").Value = 'Left'
$ws.Range("E50").Value = 'LTR'
$ws.Range("F50").Value = '0'
$ws.Range("B51").Value = 'SingleUseId78'
$ws.Range("C51").Value = 'Small'
$ws.Range("D51").Value = 'Left'
$ws.Range("E51").Value = 'LTR'
$ws.Range("F51").Value = 'Output files: <value>'
$ws.Range("B52").Value = 'SingleUseId79'
$ws.Range("C52").Value = 'Small'
$ws.Range("D52").Value = 'Left'
$ws.Range("E52").Value = 'LTR'
$ws.Range("F52").Value = '0'
